$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.ClearFormats()
}

Set-TextValue "D2" "242.83"
Set-TextValue "D4" "5.214"
Set-TextValue "D6" "3.365"
Set-TextValue "D8" "0.8055"
Set-TextValue "D9" "0.9514"
Set-TextValue "D10" "0.1427"
Set-TextValue "D11" "0.07288"
Set-TextValue "D12" "0.03127"
Set-TextValue "D13" "0.03119"
Set-TextValue "D14" "0.09277"
Set-TextValue "D15" "3.576"
Set-TextValue "D16" "0.001649"
Set-TextValue "D17" "0.04692"
Set-TextValue "D18" "0.0005799"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue "D19" "0.006354"
Set-TextValue "D20" "0.004984"
Set-TextValue "D21" "0.001042"
Set-TextValue "D23" "0.0003099"
Set-TextValue "D24" "3.759"
Set-TextValue "D25" "2.099"
$ws.Range("B28").Value = "AAXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
Set-TextValue "D28" "0.1950"
$ws.Range("E28").Value = "27AAXTokenAAB"
$ws.Range("B29").Value = "Spectre.aiUtilityToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("E29").Value = "28Spectre.aiUtilityTokenSXUT"
$ws.Range("B30").Value = "LegolasExchange"
$ws.Range("C30").Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("E30").Value = "29LegolasExchangeLGO"
$ws.Range("B31").Value = "BitZToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("E31").Value = "30BitZTokenBZ"
$ws.Range("B32").Value = "Birake"
$ws.Range("C32").Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("E32").Value = "31BirakeBIR"
$ws.Range("B33").Value = "ZBToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("E33").Value = "32ZBTokenZB"
$ws.Range("B34").Value = "NashExchange"
$ws.Range("C34").Value = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
$ws.Range("E34").Value = "33NashExchangeNEX"
Set-TextValue "D41" "0.006892"
Set-TextValue "D43" "0.1033"
Set-TextValue "D44" "0.007507"
Set-TextValue "D45" "0.00005939"
Set-TextValue "D47" "0.0005499"
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"
Set-TextValue "D48" "0.6823"
Set-TextValue "D49" "0.07378"
